# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets
# to reflect newly refreshed counts.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 26
$ws1.Range("F7").Value = 1048
$ws1.Range("F12").Value = 13373
$ws1.Range("F13").Value = 167
$ws1.Range("F14").Value = 16
$ws1.Range("F16").Value = 5498
$ws1.Range("F18").Value = 49

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 26
$ws4.Range("F29").Value = 1048
$ws4.Range("F34").Value = 13373
$ws4.Range("F35").Value = 167
$ws4.Range("F36").Value = 16
$ws4.Range("F39").Value = 5498
$ws4.Range("F41").Value = 49
